# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp banner in A1
# - Update case numbers for a handful of countries whose totals moved them
#   past a neighboring country in the (descending, by total cases) ordering,
#   which swaps which row shows which country name
# - Straight in-place numeric refresh for a few more countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Banner timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 05:26"

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Kazajistan overtakes Kuwait
Set-Row 35 "Kazajistan" 53021 1962 34149 18608 0 0 264
Set-Row 36 "Kuwait"      52007    0 42108  9520 0 0 379

# Honduras overtakes Irlanda
Set-Row 54 "Honduras" 25978 550  2721 22563 0 17 694
Set-Row 55 "Irlanda"  25542   0 23364   440 0  0 1738

# Haiti refresh (no order change)
Set-Row 85 "Haiti" 6486 54 2181 4182 0 6 123

# Mongolia refresh (no order change)
Set-Row 169 "Mongolia" 227 0 197 30 0 0 0

# Lesoto / Seychelles swap position (figures tied, only labels move)
Set-Row 184 "Lesoto"     91 0 11 80 0 0 0
Set-Row 185 "Seychelles" 91 0 11 80 0 0 0

# Antigua y Barbuda refresh (no order change)
Set-Row 189 "Antigua y Barbuda" 73 3 57 13 0 0 3

# Groenlandia / Islas Malvinas swap position (figures tied, only labels move)
Set-Row 209 "Groenlandia"    13 0 13 0 0 0 0
Set-Row 210 "Islas Malvinas" 13 0 13 0 0 0 0
